$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11557.3447621554
$ws.Range("C2").Value = 10575.4971782987
$ws.Range("D2").Value = 17616.26
$ws.Range("E2").Value = 6978.7735002457
$ws.Range("F2").Value = -2.58288839398392
$ws.Range("B3").Value = 11312.9902983077
$ws.Range("C3").Value = 10464.1358180979
$ws.Range("E3").Value = 6786.62372356818
$ws.Range("F3").Value = 250.770814236088
$ws.Range("B4").Value = 11277.8747615886
$ws.Range("C4").Value = 9814.26994216339
$ws.Range("E4").Value = 7027.59709588729
$ws.Range("F4").Value = 233.733626585445
$ws.Range("B5").Value = 4233.40721203198
$ws.Range("C5").Value = 6783.33718443316
$ws.Range("E5").Value = 6618.3288903713
$ws.Range("F5").Value = 90.391919783519
$ws.Range("B6").Value = 4134.9355805138
$ws.Range("C6").Value = 6995.9076593128
$ws.Range("E6").Value = 6682.04930302322
$ws.Range("F6").Value = 101.904040097334
$ws.Range("B7").Value = 10886.1572209584
$ws.Range("C7").Value = 10575.2172056343
$ws.Range("E7").Value = 7159.87075818216
$ws.Range("F7").Value = 270.951165159021
$ws.Range("C8").Value = 9958.12333157368
$ws.Range("F8").Value = 282.005010887713
$ws.Range("C9").Value = 9647.66447176018
$ws.Range("F9").Value = 269.06922506215
$ws.Range("C10").Value = 9400.23108100212
$ws.Range("F10").Value = 258.759500447231
$ws.Range("C11").Value = 8599.88107337639
$ws.Range("F11").Value = 225.411583462825
$ws.Range("C12").Value = 8246.57037889723
$ws.Range("F12").Value = 194.359594311261
$ws.Range("C13").Value = 8165.18089368065
$ws.Range("F13").Value = 190.618448422814
$ws.Range("C14").Value = 8826.06484264283
$ws.Range("F14").Value = 237.789015767457
$ws.Range("C15").Value = 8952.58616870612
$ws.Range("F15").Value = 243.060737686761
